$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44162
$ws.Range("J2").Value2 = 1200
$ws.Range("D3").Value2 = 44162
$ws.Range("J3").Value2 = 800
$ws.Range("K3").Value2 = 1000
$ws.Range("L3").Value2 = 1000
$ws.Range("M3").Value2 = 1000
$ws.Range("P3").Value2 = 1000
$ws.Range("D4").Value2 = 44181
$ws.Range("J4").Value2 = 1000
$ws.Range("K4").Value2 = 1300
$ws.Range("L4").Value2 = 1300
$ws.Range("M4").Value2 = 1300
$ws.Range("P4").Value2 = 1300
$ws.Range("D5").Value2 = 44181
$ws.Range("J5").Value2 = 900
$ws.Range("K5").Value2 = 900
$ws.Range("L5").Value2 = 900
$ws.Range("M5").Value2 = 900
$ws.Range("P5").Value2 = 900
$ws.Range("D6").Value2 = 44176
$ws.Range("J6").Value2 = 2500
$ws.Range("L6").Value2 = 1300
$ws.Range("M6").Value2 = 1256
$ws.Range("O6").Value2 = "Provincia de Quillota"
$ws.Range("P6").Value2 = 1256
$ws.Range("D7").Value2 = 44176
$ws.Range("I7").Value2 = "Segunda"
$ws.Range("J7").Value2 = 1500
$ws.Range("K7").Value2 = 1000
$ws.Range("L7").Value2 = 1000
$ws.Range("M7").Value2 = 1000
$ws.Range("P7").Value2 = 1000
$ws.Range("D8").Value2 = 44179
$ws.Range("I8").Value2 = "Primera"
$ws.Range("J8").Value2 = 980
$ws.Range("K8").Value2 = 1200
$ws.Range("L8").Value2 = 1200
$ws.Range("M8").Value2 = 1200
$ws.Range("O8").Value2 = "Región Metropolitana"
$ws.Range("P8").Value2 = 1200
$ws.Range("D9").Value2 = 44160
$ws.Range("J9").Value2 = 750
$ws.Range("D10").Value2 = 44160
$ws.Range("J10").Value2 = 850
$ws.Range("D11").Value2 = 44167
$ws.Range("J11").Value2 = 1430
$ws.Range("K11").Value2 = 1200
$ws.Range("M11").Value2 = 1248
$ws.Range("P11").Value2 = 1248
$ws.Range("D12").Value2 = 44167
$ws.Range("J12").Value2 = 350
$ws.Range("D13").Value2 = 44175
$ws.Range("J13").Value2 = 1500
$ws.Range("K13").Value2 = 1300
$ws.Range("M13").Value2 = 1300
$ws.Range("P13").Value2 = 1300
$ws.Range("D14").Value2 = 44175
$ws.Range("J14").Value2 = 1450
$ws.Range("D15").Value2 = 44174
$ws.Range("J15").Value2 = 2800
$ws.Range("L15").Value2 = 1250
$ws.Range("M15").Value2 = 1221
$ws.Range("P15").Value2 = 1221
$ws.Range("D16").Value2 = 44174
$ws.Range("J16").Value2 = 1300
$ws.Range("D17").Value2 = 44169
$ws.Range("J17").Value2 = 950
$ws.Range("D18").Value2 = 44169
$ws.Range("J18").Value2 = 800
$ws.Range("D21").Value2 = 44161
$ws.Range("J21").Value2 = 1600
$ws.Range("D22").Value2 = 44161
$ws.Range("J22").Value2 = 1850
$ws.Range("D23").Value2 = 44165
$ws.Range("J23").Value2 = 720
$ws.Range("K23").Value2 = 1200
$ws.Range("L23").Value2 = 1200
$ws.Range("M23").Value2 = 1200
$ws.Range("P23").Value2 = 1200
$ws.Range("D24").Value2 = 44165
$ws.Range("J24").Value2 = 750
$ws.Range("D25").Value2 = 44159
$ws.Range("J25").Value2 = 1100
$ws.Range("D26").Value2 = 44159
$ws.Range("D27").Value2 = 44172
$ws.Range("J27").Value2 = 600
$ws.Range("K27").Value2 = 1300
$ws.Range("M27").Value2 = 1300
$ws.Range("P27").Value2 = 1300
$ws.Range("D28").Value2 = 44172
$ws.Range("J28").Value2 = 550